$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H) mirrors the existing header style used by G1 ("sum")
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats - copy formatting only

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
